# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "2508" period row (row 17) entirely - shifts rows below up by one.
$ws.Rows.Item(17).Delete()

# Update dependent summary fields to reflect the remaining single period/worker row.
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1
